# Scen_NCAP_NUC_NASZ.xlsx — "Add files via upload"
#
# The real content edit in this revision is on sheet NCAP_BND: the
# Cset_CN (commodity set) column for the FLO_BND / ELE_H2 & ELE_GAS_H2
# rows (G27:G36) was renamed from "HYDROGEN" to "HYDROGEN_OUT".
# (Everything else in the shared-strings table shifts index purely as a
# side effect of that rename — we don't need to touch it by hand.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- data edit -------------------------------------------------------
$ws.Range("G27:G36").Value = "HYDROGEN_OUT"

# --- view state (scroll position / active cell at last save) --------
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("H31").Select()
